# Insert a new row at position 163 (pushes existing rows 163:265 down to 164:266)
# and populate it with the new daily price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(163).Insert()

$ws.Cells.Item(163, 1).Value = 8
$ws.Cells.Item(163, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(163, 3).Value = "Coquimbo"
$ws.Cells.Item(163, 4).Value = 44603
$ws.Cells.Item(163, 5).Value = 4
$ws.Cells.Item(163, 6).Value = 100112032
$ws.Cells.Item(163, 7).Value = "Zapallo italiano"
$ws.Cells.Item(163, 8).Value = "Sin especificar"
$ws.Cells.Item(163, 9).Value = "Primera"
$ws.Cells.Item(163, 10).Value = 520
$ws.Cells.Item(163, 11).Value = 9000
$ws.Cells.Item(163, 12).Value = 10000
$ws.Cells.Item(163, 13).Value = 9500
$ws.Cells.Item(163, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(163, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(163, 16).Value = 158
$ws.Cells.Item(163, 17).Value = 60
$ws.Cells.Item(163, 18).Value = "Hortaliza"
